# Video Recording Progress.xlsx - update progress log
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rows 2-8: G and K flip from TRUE to FALSE ---
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 7).Value = $false   # column G
    $ws.Cells.Item($r, 11).Value = $false  # column K
}

# --- Rows 14-26: I flips from FALSE to TRUE ---
for ($r = 14; $r -le 26; $r++) {
    $ws.Cells.Item($r, 9).Value = $true    # column I
}

# --- Rows 27-31: G flips to FALSE, I flips to TRUE ---
for ($r = 27; $r -le 31; $r++) {
    $ws.Cells.Item($r, 7).Value = $false   # column G
    $ws.Cells.Item($r, 9).Value = $true    # column I
}

# --- Rows 32-33: D and E flip to TRUE ---
for ($r = 32; $r -le 33; $r++) {
    $ws.Cells.Item($r, 4).Value = $true    # column D
    $ws.Cells.Item($r, 5).Value = $true    # column E
}

# --- New row 34: "Examples" section marker row ---
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = "Examples"
$ws.Cells.Item(34, 4).Value = $true
$ws.Cells.Item(34, 5).Value = $true
$ws.Cells.Item(34, 6).Value = $false
$ws.Cells.Item(34, 7).Value = $false
$ws.Cells.Item(34, 8).Value = $false
$ws.Cells.Item(34, 9).Value = $false
$ws.Cells.Item(34, 10).Value = $false
$ws.Cells.Item(34, 11).Value = $false

$hdrRange = $ws.Range("B34:C34")
$hdrRange.Font.Name = "Arial"
$hdrRange.Font.Size = 16
$hdrRange.HorizontalAlignment = -4108
$hdrRange.VerticalAlignment = -4108
$hdrRange.Merge()

$ws.Range("I9").Select()
